$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.85
$ws.Range("T2").Value = 1.71
$ws.Range("X2").Value = 19.5
$ws.Range("AO2").Value = 48
$ws.Range("F3").Value = 3.55
$ws.Range("G3").Value = 4.2
$ws.Range("H3").Value = 1.98
$ws.Range("I3").Value = 2.14
$ws.Range("J3").Value = 3.65
$ws.Range("K3").Value = 4.2
$ws.Range("N3").Value = 4.3
$ws.Range("Q3").Value = 1.71
$ws.Range("S3").Value = 2.76
$ws.Range("U3").Value = 2.26
$ws.Range("V3").Value = 1.87
$ws.Range("W3").Value = 1.31
$ws.Range("X3").Value = 19.5
$ws.Range("AA3").Value = 26
$ws.Range("AB3").Value = 18
$ws.Range("AD3").Value = 11.5
$ws.Range("AE3").Value = 21
$ws.Range("AO3").Value = 14
$ws.Range("F4").Value = 1.63
$ws.Range("K4").Value = 4.5
$ws.Range("N4").Value = 3.55
$ws.Range("O4").Value = 1.31
$ws.Range("P4").Value = 1.88
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.33
$ws.Range("S4").Value = 3.4
$ws.Range("T4").Value = 1.92
$ws.Range("U4").Value = 1.89
$ws.Range("X4").Value = 18.5
$ws.Range("AK4").Value = 980
$ws.Range("AN4").Value = 11.5
$ws.Range("H5").Value = 1.09
$ws.Range("J5").Value = 1.01
$ws.Range("Q5").Value = 1.01
$ws.Range("S5").Value = 1.05
$ws.Range("O6").Value = 1.06
$ws.Range("W6").Value = 1.58
$ws.Range("F7").Value = 4.2
$ws.Range("G7").Value = 4.9
$ws.Range("H7").Value = 1.7
$ws.Range("I7").Value = 1.84
$ws.Range("J7").Value = 4.4
$ws.Range("K7").Value = 5.2
$ws.Range("N7").Value = 5.9
$ws.Range("O7").Value = 1.16
$ws.Range("P7").Value = 2.68
$ws.Range("Q7").Value = 1.48
$ws.Range("R7").Value = 1.67
$ws.Range("S7").Value = 2.2
$ws.Range("U7").Value = 2.46
$ws.Range("V7").Value = 2.18
$ws.Range("W7").Value = 1.25
$ws.Range("X7").Value = 32
$ws.Range("AA7").Value = 21
$ws.Range("AE7").Value = 980
$ws.Range("AH7").Value = 990
$ws.Range("AJ7").Value = 110
$ws.Range("AK7").Value = 980
$ws.Range("AL7").Value = 980
$ws.Range("AN7").Value = 38
$ws.Range("AO7").Value = 7.2
$ws.Range("F8").Value = 8.199999999999999
$ws.Range("G8").Value = 9
$ws.Range("H8").Value = 1.49
$ws.Range("I8").Value = 1.52
$ws.Range("J8").Value = 4.5
$ws.Range("L8").Value = 1.42
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 3.75
$ws.Range("O8").Value = 1.34
$ws.Range("P8").Value = 1.94
$ws.Range("Q8").Value = 2.02
$ws.Range("S8").Value = 3.45
$ws.Range("T8").Value = 2.22
$ws.Range("U8").Value = 1.77
$ws.Range("V8").Value = 2.92
$ws.Range("W8").Value = 1.12
$ws.Range("X8").Value = 14.5
$ws.Range("Y8").Value = 7.4
$ws.Range("Z8").Value = 7.8
$ws.Range("AA8").Value = 13
$ws.Range("AC8").Value = 10
$ws.Range("AD8").Value = 10
$ws.Range("AH8").Value = 28
$ws.Range("AI8").Value = 48
$ws.Range("AJ8").Value = 320
$ws.Range("AK8").Value = 160
$ws.Range("AL8").Value = 150
$ws.Range("AM8").Value = 190
$ws.Range("AN8").Value = 240
$ws.Range("AO8").Value = 9
$ws.Range("F9").Value = 1.66
$ws.Range("G9").Value = 1.73
$ws.Range("J9").Value = 3.65
$ws.Range("K9").Value = 3.95
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 2.96
$ws.Range("O9").Value = 1.46
$ws.Range("P9").Value = 1.65
$ws.Range("Q9").Value = 2.34
$ws.Range("R9").Value = 1.24
$ws.Range("S9").Value = 4.5
$ws.Range("T9").Value = 2.2
$ws.Range("U9").Value = 1.71
$ws.Range("V9").Value = 1.16
$ws.Range("W9").Value = 2.36
$ws.Range("X9").Value = 13
$ws.Range("Z9").Value = 55
$ws.Range("AA9").Value = 300
$ws.Range("AB9").Value = 6.6
$ws.Range("AC9").Value = 8.6
$ws.Range("AE9").Value = 160
$ws.Range("AF9").Value = 9
$ws.Range("AI9").Value = 170
$ws.Range("AJ9").Value = 980
$ws.Range("AL9").Value = 55
$ws.Range("AM9").Value = 260
$ws.Range("AN9").Value = 1000
$ws.Range("F11").Value = 6.8
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = 1.59
$ws.Range("I11").Value = 1.6
$ws.Range("V11").Value = 2.66
$ws.Range("AA11").Value = 14.5
$ws.Range("AO11").Value = 9.4
$ws.Range("G12").Value = 2.26
$ws.Range("I12").Value = 4.6
$ws.Range("M12").Value = 1.17
$ws.Range("Q12").Value = 3.35
$ws.Range("R12").Value = 1.13
$ws.Range("S12").Value = 7.8
$ws.Range("W12").Value = 1.79
$ws.Range("F13").Value = 1.85
$ws.Range("G13").Value = 1.88
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 5.4
$ws.Range("J13").Value = 3.65
$ws.Range("K13").Value = 3.75
$ws.Range("N13").Value = 3.3
$ws.Range("P13").Value = 1.8
$ws.Range("R13").Value = 1.3
$ws.Range("S13").Value = 3.9
$ws.Range("T13").Value = 1.97
$ws.Range("V13").Value = 1.22
$ws.Range("W13").Value = 2.12
$ws.Range("AA13").Value = 150
$ws.Range("AE13").Value = 85
$ws.Range("AF13").Value = 10.5
$ws.Range("AJ13").Value = 20
$ws.Range("AK13").Value = 21
$ws.Range("AM13").Value = 150
$ws.Range("AN13").Value = 15
